# Slide 5, shape 1 ("Google Shape;144;p20") is the body placeholder whose
# first bullet currently reads:
#   "Matched 119,656 bibliographic records."
# The commit extends that sentence (dropping the trailing period from the
# first run and re-adding it as its own run at the very end) to:
#   "Matched 119,656 bibliographic records in two json files (which you can
#    find on GitHub)."

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$para1 = $tr.Paragraphs(1, 1)

# --- 1. Strip the trailing "." off the existing run -----------------------
# Paragraphs()/TextRange.Text reports the paragraph mark as part of the
# string, so trim CR/LF/VT before measuring the real text length.
$rawText = $para1.Text
$coreText = $rawText.TrimEnd([char]13, [char]10, [char]11)
if ($coreText.EndsWith(".")) {
    $periodPos = $para1.Start + $coreText.Length - 1
    $tr.Characters($periodPos, 1).Text = ""
}

# --- 2. Re-fetch paragraph 1 (its length changed) and append the new runs -
$para1 = $tr.Paragraphs(1, 1)

$para1.InsertAfter(" in ")
$para1.InsertAfter("two")
$para1.InsertAfter(" ")
$para1.InsertAfter("json")
$para1.InsertAfter(" ")
$para1.InsertAfter("files")
$para1.InsertAfter(" (")
$para1.InsertAfter("which")
$para1.InsertAfter(" ")
$para1.InsertAfter("you")
$para1.InsertAfter(" ")
$para1.InsertAfter("can")
$para1.InsertAfter(" ")
$para1.InsertAfter("find")
$para1.InsertAfter(" ")
$para1.InsertAfter("on GitHub)")
$para1.InsertAfter(".")
